$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "{2.0: 3013, 1.0: 2987}"
$ws.Range("C2").Value = "{2.0: 50.22, 1.0: 49.78}"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2.0"

$ws.Range("B3").Value = "{1.0: 3012, 2.0: 2988}"
$ws.Range("C3").Value = "{1.0: 50.2, 2.0: 49.8}"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.0"

$ws.Range("B4").Value = "{1.0: 3034, 2.0: 2966}"
$ws.Range("C4").Value = "{1.0: 50.57, 2.0: 49.43}"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0"

$ws.Range("B5").Value = "{1.0: 2011, 2.0: 1998, 3.0: 1991}"
$ws.Range("C5").Value = "{1.0: 33.52, 2.0: 33.3, 3.0: 33.18}"

$ws.Range("B6").Value = "{2.0: 3003, 1.0: 2997}"
$ws.Range("C6").Value = "{2.0: 50.05, 1.0: 49.95}"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.0"

$ws.Range("B7").Value = "{2.0: 2038, 1.0: 2010, 3.0: 1952}"
$ws.Range("C7").Value = "{2.0: 33.97, 1.0: 33.5, 3.0: 32.53}"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.0"

$ws.Range("B8").Value = "{2.0: 3063, 1.0: 2937}"
$ws.Range("C8").Value = "{2.0: 51.05, 1.0: 48.95}"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.0"

$ws.Range("B9").Value = "{2.0: 3029, 1.0: 2971}"
$ws.Range("C9").Value = "{2.0: 50.48, 1.0: 49.52}"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.0"

$ws.Range("B10").Value = "{1.0: 3006, 2.0: 2994}"
$ws.Range("C10").Value = "{1.0: 50.1, 2.0: 49.9}"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.0"

$ws.Range("B11").Value = "{2.0: 3098, 1.0: 2902}"
$ws.Range("C11").Value = "{2.0: 51.63, 1.0: 48.37}"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.0"

$ws.Range("B12").Value = "{1.0: 3084, 2.0: 2916}"
$ws.Range("C12").Value = "{1.0: 51.4, 2.0: 48.6}"

$ws.Range("B13").Value = "{1.0: 3006, 2.0: 2994}"
$ws.Range("C13").Value = "{1.0: 50.1, 2.0: 49.9}"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.0"

$ws.Range("B14").Value = "{2.0: 3045, 1.0: 2955}"
$ws.Range("C14").Value = "{2.0: 50.75, 1.0: 49.25}"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.0"

$ws.Range("B15").Value = "{2.0: 3009, 1.0: 2991}"
$ws.Range("C15").Value = "{2.0: 50.15, 1.0: 49.85}"

$ws.Range("B16").Value = "{2.0: 3035, 1.0: 2965}"
$ws.Range("C16").Value = "{2.0: 50.58, 1.0: 49.42}"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.0"

$ws.Range("B17").Value = "{2.0: 3056, 1.0: 2944}"
$ws.Range("C17").Value = "{2.0: 50.93, 1.0: 49.07}"

$ws.Range("B18").Value = "{2.0: 3070, 1.0: 2930}"
$ws.Range("C18").Value = "{2.0: 51.17, 1.0: 48.83}"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.0"

$ws.Range("B19").Value = "{2.0: 3070, 1.0: 2930}"
$ws.Range("C19").Value = "{2.0: 51.17, 1.0: 48.83}"

$ws.Range("B20").Value = "{3.0: 2022, 1.0: 1989, 2.0: 1989}"
$ws.Range("C20").Value = "{3.0: 33.7, 1.0: 33.15, 2.0: 33.15}"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.0"

$ws.Range("B21").Value = "{1.0: 2021, 3.0: 2007, 2.0: 1972}"
$ws.Range("C21").Value = "{1.0: 33.68, 3.0: 33.45, 2.0: 32.87}"
